# Fruta / hortaliza, semanal
# Insert a new weekly data row for "Vega Modelo de Temuco - Mango" right
# before the current row 454, pushing the existing rows 454:518 down to
# 455:519 (dimension grows from A1:T518 to A1:T519).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 454 (shifts 454:518 -> 455:519).
$ws.Rows.Item(454).Insert()

# Populate the newly inserted row 454 with the new weekly record.
$ws.Cells.Item(454, 1).Value = 10
$ws.Cells.Item(454, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(454, 3).Value = 'La Araucanía'
$ws.Cells.Item(454, 4).Value = 44984
$ws.Cells.Item(454, 5).Value = 9
$ws.Cells.Item(454, 6).Value = 'Fruta'
$ws.Cells.Item(454, 7).Value = 100108
$ws.Cells.Item(454, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(454, 9).Value = 100108002
$ws.Cells.Item(454, 10).Value = 'Mango'
$ws.Cells.Item(454, 11).Value = 'Sin especificar'
$ws.Cells.Item(454, 12).Value = 'Primera'
$ws.Cells.Item(454, 13).Value = 125
$ws.Cells.Item(454, 14).Value = 8500
$ws.Cells.Item(454, 15).Value = 8500
$ws.Cells.Item(454, 16).Value = 8500
$ws.Cells.Item(454, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(454, 18).Value = 'Perú'
$ws.Cells.Item(454, 19).Value = 2125
$ws.Cells.Item(454, 20).Value = 4
